# LNF LNA Shipments - add a new shipment batch (#7) to the log on Sheet1.
#
# The sheet tracks shipment "sections" that each start with a bold header
# row (section #, date received, first item) followed by continuation rows
# for each additional item in the batch, and end in a thick-bottom-border
# blank row. The last existing section (#6) occupies rows 43-53. We copy
# that section's row layout/formatting down to build section #7 (rows
# 54-60), then overwrite the cells with the new batch's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clone the formatting of section 6's header + 5 continuation rows (43-48)
# into the new section's equivalent rows (54-59), and clone the closing
# blank/thick-border row (53) into the new closing row (60).
$ws.Range("B43:G48").Copy($ws.Range("B54"))
$ws.Range("B53:G53").Copy($ws.Range("B60"))
$excel.CutCopyMode = 0

# Header row of the new section: batch number, date received, first qty.
$ws.Range("B54").Value = 7
$ws.Range("C54").Value = "2021-03-25"

# LNF serial numbers (column E) and Minex serial numbers (column F) for
# the six units in this batch - fill column E fully, then column F fully,
# so new shared-string entries are interleaved the same way Excel would
# add them (LNF strings first, then Minex strings).
$lnfSerials = @(
    "LNF-ABLNC1_15A sn0027A_v2",
    "LNF-ABLNC1_15A sn0032A_v2",
    "LNF-ABLNC1_15A sn0058A_v2",
    "LNF-ABLNC1_15A sn0066A_v2",
    "LNF-ABLNC1_15A sn0067A_v2",
    "LNF-ABLNC1_15A sn0098A_v2"
)
$minexSerials = @(
    "C-0027F",
    "C-0032F",
    "C-0058F",
    "C-0066F",
    "C-0067F",
    "C-0098F"
)

for ($i = 0; $i -lt $lnfSerials.Length; $i++) {
    $row = 54 + $i
    $ws.Range("E$row").Value = $lnfSerials[$i]
}
for ($i = 0; $i -lt $minexSerials.Length; $i++) {
    $row = 54 + $i
    $ws.Range("F$row").Value = $minexSerials[$i]
}

# Match the saved selection/view state left behind after the edit.
$ws.Range("B2:G60").Select()
